$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original "email" cell font so it can be restored after the
# hyperlinks are rebuilt (Hyperlinks.Add forces Excel's built-in Hyperlink
# theme style onto the target cell).
$emailFontColor = $ws.Range("I2").Font.Color
$emailFontUnderline = $ws.Range("I2").Font.Underline
$emailFontName = $ws.Range("I2").Font.Name
$emailFontSize = $ws.Range("I2").Font.Size

# Remove the existing hyperlinks before the columns shift, since the engine
# does not re-target a hyperlink's range when EntireColumn.Delete() shifts
# cell content underneath it.
$ws.Range("I2").Hyperlinks.Delete()

# Delete the whole "jezyk" column (F). Everything to the right (grupa, nr
# tel, email, notatka rekrutacyjna) shifts one column to the left.
$ws.Range("F1").EntireColumn.Delete()

# Re-create the two mailto hyperlinks at their new location (now column H).
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:olo@gmail.com", [Type]::Missing, [Type]::Missing, "olo@gmail.com")
$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:ddfdfd@fgdg", [Type]::Missing, [Type]::Missing, "ddfdfd@fgdg")

# Restore the original (non-themed) hyperlink cell formatting.
foreach ($addr in @("H2", "H3")) {
    $c = $ws.Range($addr)
    $c.Font.Color = $emailFontColor
    $c.Font.Underline = $emailFontUnderline
    $c.Font.Name = $emailFontName
    $c.Font.Size = $emailFontSize
}

$ws.Range("F1").Select()
